$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "296.99"
Set-TextValue "E2" "1.62%"
Set-TextValue "D3" "41.83"
Set-TextValue "E3" "3.76%"
Set-TextValue "D4" "5.000"
Set-TextValue "E4" "-0.22%"
Set-TextValue "D5" "0.07514"
Set-TextValue "E5" "2.68%"
Set-TextValue "D6" "4.375"
Set-TextValue "E6" "1.75%"
Set-TextValue "D7" "1.583"
Set-TextValue "E7" "3.46%"
Set-TextValue "D8" "0.9257"
Set-TextValue "E8" "-0.14%"
Set-TextValue "E9" "1.36%"
Set-TextValue "D10" "0.1194"
Set-TextValue "E10" "0.65%"
Set-TextValue "D11" "0.1821"
Set-TextValue "E11" "4.44%"
Set-TextValue "D12" "0.08904"
Set-TextValue "E12" "2.40%"
Set-TextValue "D13" "0.04078"
Set-TextValue "E13" "-5.69%"
Set-TextValue "D14" "0.1049"
Set-TextValue "E14" "-0.52%"
Set-TextValue "D15" "0.001277"
Set-TextValue "E15" "0.97%"
Set-TextValue "D16" "0.005780"
Set-TextValue "E16" "-3.74%"
Set-TextValue "D17" "3.355"
Set-TextValue "E17" "0.48%"
Set-TextValue "E18" "0.76%"
Set-TextValue "D19" "8.094"
Set-TextValue "E19" "1.53%"
Set-TextValue "D20" "0.1391"
Set-TextValue "E20" "0.03%"
Set-TextValue "E21" "11.06%"
Set-TextValue "D22" "0.04103"
Set-TextValue "E22" "4.37%"
Set-TextValue "D23" "0.001268"
Set-TextValue "E23" "0.64%"
Set-TextValue "E24" "3.14%"
Set-TextValue "E25" "-3.93%"
Set-TextValue "D38" "0.02403"
Set-TextValue "E38" "5.56%"
Set-TextValue "D39" "0.05200"
Set-TextValue "E39" "4.48%"
Set-TextValue "D41" "0.007799"
Set-TextValue "E41" "1.07%"
Set-TextValue "E42" "3.23%"
Set-TextValue "D43" "0.007411"
Set-TextValue "E43" "0.57%"
Set-TextValue "E44" "-0.44%"
Set-TextValue "E45" "1.49%"
Set-TextValue "D46" "0.00006600"
Set-TextValue "E46" "4.77%"
Set-TextValue "E47" "-0.02%"
Set-TextValue "D48" "0.03163"
Set-TextValue "E48" "48.33%"
Set-TextValue "E49" "0.06%"
Set-TextValue "E50" "-0.02%"
Set-TextValue "E51" "-0.02%"
